{"js": "// Replace the 25 \"NNN\u00d7N=\" equation texts in the practice-sheet table with\n// their new values, one-for-one, in document order. Each old value is a\n// unique string in the document, so a targeted search + in-place replace\n// keeps every run's original formatting (font / size) untouched.\nconst replacements = [\n  [\"970\u00d74=\", \"825\u00d78=\"],\n  [\"493\u00d77=\", \"559\u00d78=\"],\n  [\"171\u00d72=\", \"489\u00d77=\"],\n  [\"367\u00d73=\", \"364\u00d77=\"],\n  [\"572\u00d77=\", \"187\u00d78=\"],\n  [\"987\u00d74=\", \"575\u00d76=\"],\n  [\"981\u00d73=\", \"583\u00d75=\"],\n  [\"117\u00d79=\", \"273\u00d73=\"],\n  [\"588\u00d75=\", \"777\u00d73=\"],\n  [\"218\u00d75=\", \"539\u00d74=\"],\n  [\"471\u00d77=\", \"870\u00d72=\"],\n  [\"491\u00d78=\", \"875\u00d79=\"],\n  [\"548\u00d73=\", \"382\u00d73=\"],\n  [\"210\u00d77=\", \"847\u00d75=\"],\n  [\"938\u00d75=\", \"893\u00d75=\"],\n  [\"527\u00d73=\", \"861\u00d76=\"],\n  [\"693\u00d74=\", \"654\u00d75=\"],\n  [\"276\u00d78=\", \"315\u00d78=\"],\n  [\"187\u00d72=\", \"741\u00d72=\"],\n  [\"336\u00d79=\", \"528\u00d75=\"],\n  [\"610\u00d76=\", \"768\u00d78=\"],\n  [\"523\u00d74=\", \"972\u00d72=\"],\n  [\"379\u00d73=\", \"387\u00d76=\"],\n  [\"712\u00d79=\", \"160\u00d78=\"],\n  [\"596\u00d78=\", \"367\u00d74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"NNN\u00d7N=\" equation texts in the practice-sheet table with\n# their new values, one-for-one, in document order. Each old value is a\n# unique string in the document, so Find/Replace targeted at that exact\n# string only ever touches the single matching run and leaves its\n# formatting (font / size) untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"970\u00d74=\", \"825\u00d78=\"),\n    @(\"493\u00d77=\", \"559\u00d78=\"),\n    @(\"171\u00d72=\", \"489\u00d77=\"),\n    @(\"367\u00d73=\", \"364\u00d77=\"),\n    @(\"572\u00d77=\", \"187\u00d78=\"),\n    @(\"987\u00d74=\", \"575\u00d76=\"),\n    @(\"981\u00d73=\", \"583\u00d75=\"),\n    @(\"117\u00d79=\", \"273\u00d73=\"),\n    @(\"588\u00d75=\", \"777\u00d73=\"),\n    @(\"218\u00d75=\", \"539\u00d74=\"),\n    @(\"471\u00d77=\", \"870\u00d72=\"),\n    @(\"491\u00d78=\", \"875\u00d79=\"),\n    @(\"548\u00d73=\", \"382\u00d73=\"),\n    @(\"210\u00d77=\", \"847\u00d75=\"),\n    @(\"938\u00d75=\", \"893\u00d75=\"),\n    @(\"527\u00d73=\", \"861\u00d76=\"),\n    @(\"693\u00d74=\", \"654\u00d75=\"),\n    @(\"276\u00d78=\", \"315\u00d78=\"),\n    @(\"187\u00d72=\", \"741\u00d72=\"),\n    @(\"336\u00d79=\", \"528\u00d75=\"),\n    @(\"610\u00d76=\", \"768\u00d78=\"),\n    @(\"523\u00d74=\", \"972\u00d72=\"),\n    @(\"379\u00d73=\", \"387\u00d76=\"),\n    @(\"712\u00d79=\", \"160\u00d78=\"),\n    @(\"596\u00d78=\", \"367\u00d74=\")\n)\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n}\n"}
